$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")
$ws.Range("A4").Value = "Y"
$ws.Range("A8").Value = "Y"
$ws.Range("A11").Value = "Y"
